# Updated cryptos list (Price/Volume(1h) refresh) on Sat Oct  5 08:58:58 UTC 2024
# with GitHub Actions.
#
# Cells that hold plain-numeric-looking text (e.g. "564.00") are first
# switched to Text format ("@") so Excel keeps them as literal strings
# instead of auto-converting to numbers (matching the source workbook,
# where every Price/Volume cell is an inline/shared string, not a number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.169.12'
$ws.Range("E2").Value = '  +1.07%  '
$ws.Range("D3").Value = '2.421.15'
$ws.Range("E3").Value = '  +1.69%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '564.00'
$ws.Range("E5").Value = '  +2.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.45'
$ws.Range("E6").Value = '  +3.04%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  +1.49%  '
$ws.Range("D9").Value = '2.418.57'
$ws.Range("E9").Value = '  +1.52%  '
$ws.Range("E10").Value = '  +1.10%  '
$ws.Range("E11").Value = '  -1.67%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.41'
$ws.Range("E12").Value = '  +1.14%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.355'
$ws.Range("E13").Value = '  +1.22%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.05'
$ws.Range("E14").Value = '  +2.25%  '
$ws.Range("E15").Value = '  +5.18%  '
$ws.Range("D16").Value = '2.858.35'
$ws.Range("E16").Value = '  +2.30%  '
$ws.Range("D17").Value = '61.984.12'
$ws.Range("E17").Value = '  +0.92%  '
$ws.Range("D18").Value = '2.418.46'
$ws.Range("E18").Value = '  +1.58%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.29'
$ws.Range("E19").Value = '  +2.69%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '324.94'
$ws.Range("E20").Value = '  +1.34%  '
$ws.Range("E21").Value = '  +1.18%  '
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("E23").Value = '  +0.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.49'
$ws.Range("E24").Value = '  +1.63%  '
$ws.Range("E25").Value = '  -1.69%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.06'
$ws.Range("E26").Value = '  +1.66%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '586.70'
$ws.Range("E27").Value = '  +12.57%  '
$ws.Range("E28").Value = '  +1.74%  '
$ws.Range("E29").Value = '  +0.21%  '
$ws.Range("D30").Value = '0.0₃0948'
$ws.Range("E30").Value = '  +4.72%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.25'
$ws.Range("E31").Value = '  +0.40%  '
$ws.Range("E32").Value = '  +4.67%  '
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("E34").Value = '  +2.27%  '
$ws.Range("E35").Value = '  +0.59%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.73'
$ws.Range("E36").Value = '  +2.85%  '
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("E38").Value = '  +1.50%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '154.49'
$ws.Range("E39").Value = '  +5.35%  '
$ws.Range("E40").Value = '  +1.21%  '
$ws.Range("E41").Value = '  +0.90%  '
$ws.Range("E42").Value = '  -3.61%  '
$ws.Range("E43").Value = '  +0.07%  '
$ws.Range("E44").Value = '  +8.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '149.99'
$ws.Range("E45").Value = '  +0.95%  '
$ws.Range("E46").Value = '  +1.16%  '
$ws.Range("E47").Value = '  +2.38%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.41'
$ws.Range("E48").Value = '  +3.41%  '
$ws.Range("E49").Value = '  +1.58%  '
$ws.Range("E50").Value = '  +2.04%  '
$ws.Range("E51").Value = '  +1.83%  '
